$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) B6: new inline text "a", bold + centered (matches the new font/style introduced
#    for the auto-complete-looking helper cell while keeping the existing border).
$b6 = $ws.Range("B6")
$b6.Value = "a"
$b6.Font.Bold = $true
$b6.HorizontalAlignment = -4108   # xlCenter
$b6.VerticalAlignment = -4108     # xlCenter

# 2) Opening Comments / Concluding Comments cells: drop the trailing counter digit,
#    leaving a lone space after the line break (colour eye-dropper clean-up).
$ws.Range("C6").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C21").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C28").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C43").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C50").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C66").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C73").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C89").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C98").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C114").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C121").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C138").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C145").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C162").Value = "Concluding Comments " + [char]10 + " "
$ws.Range("C169").Value = "Opening Comments " + [char]10 + " "
$ws.Range("C185").Value = "Concluding Comments " + [char]10 + " "

# 3) C60 / C83 ("Bible Study" slots): match the wrap-text formatting used by the
#    neighbouring Initial Call / Return Visit rows in the same weeks.
$ws.Range("C60").WrapText = $true
$ws.Range("C83").WrapText = $true
